$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$ws.Range("I31").Value = "best+variance sacling+regularization"
$ws.Range("I32").Value = "Epoch 22 trained in 416.028952, cost 113.340152"
$ws.Range("I33").Value = "Model saved in file: ./3d-object-recognition\ShapeNet\ShapeNet22.ckpt"
$ws.Range("I34").Value = "train deconvolution average accuracy 0.957705"
$ws.Range("I35").Value = "train category accuracy 0.993248"
$ws.Range("I36").Value = "Evaluating train"
$ws.Range("I37").Value = "Category airplane has 4 parts and IoU 0.826747"
$ws.Range("I38").Value = "Category bag has 2 parts and IoU 0.808461"
$ws.Range("I39").Value = "Category cap has 2 parts and IoU 0.916903"
$ws.Range("I40").Value = "Category car has 4 parts and IoU 0.791826"
$ws.Range("I41").Value = "Category chair has 4 parts and IoU 0.911943"
$ws.Range("I42").Value = "Category earphone has 3 parts and IoU 0.556492"
$ws.Range("I43").Value = "Category guitar has 3 parts and IoU 0.881702"
$ws.Range("I44").Value = "Category knife has 2 parts and IoU 0.844052"
$ws.Range("I45").Value = "Category lamp has 4 parts and IoU 0.870830"
$ws.Range("I46").Value = "Category laptop has 2 parts and IoU 0.956305"
$ws.Range("I47").Value = "Category motorbike has 6 parts and IoU 0.556993"
$ws.Range("I48").Value = "Category mug has 2 parts and IoU 0.937098"
$ws.Range("I49").Value = "Category pistol has 3 parts and IoU 0.851671"
$ws.Range("I50").Value = "Category rocket has 3 parts and IoU 0.441905"
$ws.Range("I51").Value = "Category skateboard has 3 parts and IoU 0.730330"
$ws.Range("I52").Value = "Category table has 3 parts and IoU 0.861034"
$ws.Range("I53").Value = "Weighted average IOU is 0.860397"
$ws.Range("I54").Value = "dev deconvolution average accuracy 0.939193"
$ws.Range("I55").Value = "dev category accuracy 0.989774"
$ws.Range("I56").Value = "Evaluating dev"
$ws.Range("I57").Value = "Category airplane has 4 parts and IoU 0.833135"
$ws.Range("I58").Value = "Category bag has 2 parts and IoU 0.739287"
$ws.Range("I59").Value = "Category cap has 2 parts and IoU 0.542882"
$ws.Range("I60").Value = "Category car has 4 parts and IoU 0.777862"
$ws.Range("I61").Value = "Category chair has 4 parts and IoU 0.903612"
$ws.Range("I62").Value = "Category earphone has 3 parts and IoU 0.511505"
$ws.Range("I63").Value = "Category guitar has 3 parts and IoU 0.879616"
$ws.Range("I64").Value = "Category knife has 2 parts and IoU 0.781302"
$ws.Range("I65").Value = "Category lamp has 4 parts and IoU 0.796643"
$ws.Range("I66").Value = "Category laptop has 2 parts and IoU 0.965612"
$ws.Range("I67").Value = "Category motorbike has 6 parts and IoU 0.525037"
$ws.Range("I68").Value = "Category mug has 2 parts and IoU 0.925973"
$ws.Range("I69").Value = "Category pistol has 3 parts and IoU 0.841207"
$ws.Range("I70").Value = "Category rocket has 3 parts and IoU 0.437171"
$ws.Range("I71").Value = "Category skateboard has 3 parts and IoU 0.631756"
$ws.Range("I72").Value = "Category table has 3 parts and IoU 0.836823"
$ws.Range("I73").Value = "Weighted average IOU is 0.839242"

# Update the view: scroll so row 19 is the top-left visible row, and select I32
# (matches topLeftCell="A19" / activeCell="I32" sqref="I32" in the target sheetView)
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 19
    $win.ScrollColumn = 1
} catch {
}
$ws.Range("I32").Select()
